# EPBDS-3135 Check for return type in conditional array index.
# Adds a new "errorSelect" / "errorSelectLiteral" example (rows 21-22),
# mirroring the existing driverSelectMany / driverSelectManyLiteral block
# (rows 17-18) in both layout and formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the formatting of the existing two-row example block (rows 17:18)
# for the two new rows (21:22) - copy one row at a time so the engine
# reuses the existing style/font table entries instead of minting new
# (visually identical) duplicates.
$ws.Range("B17:D17").Copy() | Out-Null
$ws.Range("B21:D21").PasteSpecial(-4122) | Out-Null
$ws.Range("B18:D18").Copy() | Out-Null
$ws.Range("B22:D22").PasteSpecial(-4122) | Out-Null

$ws.Range("F17:H17").Copy() | Out-Null
$ws.Range("F21:H21").PasteSpecial(-4122) | Out-Null
$ws.Range("F18:H18").Copy() | Out-Null
$ws.Range("F22:H22").PasteSpecial(-4122) | Out-Null

# New example text content.
$ws.Range("B21").Value = "Method Driver[] errorSelect(Driver[] arrayOfDrivers)"
$ws.Range("B22").Value = "return arrayOfDrivers[@ age = 20];"
$ws.Range("F21").Value = "Method Driver[] errorSelectLiteral(Driver[] arrayOfDrivers)"
$ws.Range("F22").Value = "return arrayOfDrivers[select all having  numMovingViolations = 0];"

# Match row 18's explicit row height on the new second row (row 22).
$ws.Rows.Item(22).RowHeight = 15

# Merge the label/description cells, same pattern as the other blocks.
$ws.Range("B21:D21").Merge() | Out-Null
$ws.Range("B22:D22").Merge() | Out-Null
$ws.Range("F21:H21").Merge() | Out-Null
$ws.Range("F22:H22").Merge() | Out-Null
